$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-5: runs/balls/fours/sixes (columns C:F) get updated -- the
# activity figures for row 2 and row 5 are swapped, and likewise for
# row 3 and row 4.
#
# The leading apostrophe forces Excel to keep the value stored as text
# (matching the sheet's existing "number stored as text" cells) instead
# of silently converting it to a numeric cell; re-applying the "Normal"
# style afterwards clears the quote-prefix formatting flag that the
# apostrophe trick adds, so the cell format stays exactly as it was.

function Set-TextValue($rng, $val) {
    $rng.Value = "'" + $val
    $rng.Style = "Normal"
}

Set-TextValue $ws.Range("C2") "12"
Set-TextValue $ws.Range("D2") "9"
Set-TextValue $ws.Range("E2") "0"
Set-TextValue $ws.Range("F2") "1"

Set-TextValue $ws.Range("C3") "7"
Set-TextValue $ws.Range("D3") "5"
Set-TextValue $ws.Range("E3") "1"
Set-TextValue $ws.Range("F3") "0"

Set-TextValue $ws.Range("C4") "14"
Set-TextValue $ws.Range("D4") "13"
Set-TextValue $ws.Range("E4") "2"
Set-TextValue $ws.Range("F4") "0"

Set-TextValue $ws.Range("C5") "0"
Set-TextValue $ws.Range("D5") "0"
Set-TextValue $ws.Range("E5") "0"
Set-TextValue $ws.Range("F5") "0"
